# Auto-generated Excel COM-interop edit script
# Applies cell value changes per sheet as described by the authoritative diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 540.7143  # H2: 607.5 -> 540.7143
$ws.Cells.Item(2, 10).Value = 534.75  # J2: 666.3333 -> 534.75
$ws.Cells.Item(2, 12).Value = 534.75  # L2: 666.3333 -> 534.75
$ws.Cells.Item(2, 14).Value = -760.75  # N2: -892.3333 -> -760.75
$ws.Cells.Item(18, 8).Value = 324.75  # H18: 339.8 -> 324.75
$ws.Cells.Item(18, 9).Value = 324.75  # I18: 339.8 -> 324.75
$ws.Cells.Item(18, 11).Value = 324.75  # K18: 339.8 -> 324.75
$ws.Cells.Item(18, 13).Value = -40.75  # M18: -55.80000000000001 -> -40.75
$ws.Cells.Item(21, 8).Value = 0  # H21: 2500 -> 0
$ws.Cells.Item(21, 9).Value = 0  # I21: 2500 -> 0
$ws.Cells.Item(21, 11).Value = 0  # K21: 2500 -> 0
$ws.Cells.Item(21, 13).ClearContents()  # M21: -2032 -> (cleared)
$ws.Cells.Item(23, 8).Value = 0  # H23: 2500 -> 0
$ws.Cells.Item(23, 9).Value = 0  # I23: 2500 -> 0
$ws.Cells.Item(23, 11).Value = 0  # K23: 2500 -> 0
$ws.Cells.Item(23, 13).ClearContents()  # M23: -2266 -> (cleared)
$ws.Cells.Item(29, 8).Value = 7503.25  # H29: 5.3333335 -> 7503.25
$ws.Cells.Item(29, 9).Value = 6.5  # I29: 5.3333335 -> 6.5
$ws.Cells.Item(29, 10).Value = 15000  # J29: 0 -> 15000
$ws.Cells.Item(29, 11).Value = 19.5  # K29: 16.0000005 -> 19.5
$ws.Cells.Item(29, 12).Value = 45000  # L29: 0 -> 45000
$ws.Cells.Item(29, 13).Value = 261.5  # M29: 264.9999995 -> 261.5
$ws.Cells.Item(29, 14).Value = -45562  # N29: (empty) -> -45562
$ws.Cells.Item(38, 8).Value = 967.94116  # H38: 1483.4286 -> 967.94116
$ws.Cells.Item(38, 9).Value = 163.66667  # I38: 591.7778 -> 163.66667
$ws.Cells.Item(38, 10).Value = 7000  # J38: 6833.3335 -> 7000
$ws.Cells.Item(38, 11).Value = 491.00001  # K38: 1775.3334 -> 491.00001
$ws.Cells.Item(38, 12).Value = 21000  # L38: 20500.0005 -> 21000
$ws.Cells.Item(38, 13).Value = -119.00001  # M38: -1403.3334 -> -119.00001
$ws.Cells.Item(38, 14).Value = -21744  # N38: -21244.0005 -> -21744
$ws.Cells.Item(40, 8).Value = 3659.6667  # H40: 3729.6453 -> 3659.6667
$ws.Cells.Item(40, 9).Value = 3037  # I40: 3000 -> 3037
$ws.Cells.Item(40, 10).Value = 4033.2666  # J40: 3983.4348 -> 4033.2666
$ws.Cells.Item(40, 11).Value = 3037  # K40: 3000 -> 3037
$ws.Cells.Item(40, 12).Value = 4033.2666  # L40: 3983.4348 -> 4033.2666
$ws.Cells.Item(40, 13).Value = -2862  # M40: -2825 -> -2862
$ws.Cells.Item(40, 14).Value = -4383.2666  # N40: -4333.4348 -> -4383.2666
$ws.Cells.Item(43, 8).Value = 3999.6  # H43: 4499.8 -> 3999.6
$ws.Cells.Item(51, 8).Value = 13015.3  # H51: 11266.083 -> 13015.3
$ws.Cells.Item(51, 9).Value = 15281.5  # I51: 13782.444 -> 15281.5
$ws.Cells.Item(51, 10).Value = 3950.5  # J51: 3717 -> 3950.5
$ws.Cells.Item(51, 11).Value = 15281.5  # K51: 13782.444 -> 15281.5
$ws.Cells.Item(51, 12).Value = 3950.5  # L51: 3717 -> 3950.5
$ws.Cells.Item(51, 13).Value = -14797.5  # M51: -13298.444 -> -14797.5
$ws.Cells.Item(51, 14).Value = -4918.5  # N51: -4685 -> -4918.5
$ws.Cells.Item(76, 8).Value = 5039.8  # H76: 5221.8887 -> 5039.8
$ws.Cells.Item(76, 10).Value = 5099.75  # J76: 5666 -> 5099.75
$ws.Cells.Item(76, 12).Value = 5099.75  # L76: 5666 -> 5099.75
$ws.Cells.Item(76, 14).Value = -5729.75  # N76: -6296 -> -5729.75
$ws.Cells.Item(79, 8).Value = 5039.8  # H79: 5221.8887 -> 5039.8
$ws.Cells.Item(79, 10).Value = 5099.75  # J79: 5666 -> 5099.75
$ws.Cells.Item(79, 12).Value = 5099.75  # L79: 5666 -> 5099.75
$ws.Cells.Item(79, 14).Value = -7283.75  # N79: -7850 -> -7283.75
$ws.Cells.Item(86, 8).Value = 2755.3157  # H86: 2881.35 -> 2755.3157
$ws.Cells.Item(86, 9).Value = 2433.8333  # I86: 2643.1428 -> 2433.8333
$ws.Cells.Item(86, 10).Value = 2903.6924  # J86: 3009.6155 -> 2903.6924
$ws.Cells.Item(86, 11).Value = 2433.8333  # K86: 2643.1428 -> 2433.8333
$ws.Cells.Item(86, 12).Value = 2903.6924  # L86: 3009.6155 -> 2903.6924
$ws.Cells.Item(86, 13).Value = -1310.8333  # M86: -1520.1428 -> -1310.8333
$ws.Cells.Item(86, 14).Value = -5149.6924  # N86: -5255.6155 -> -5149.6924
$ws.Cells.Item(89, 8).Value = 2755.3157  # H89: 2881.35 -> 2755.3157
$ws.Cells.Item(89, 9).Value = 2433.8333  # I89: 2643.1428 -> 2433.8333
$ws.Cells.Item(89, 10).Value = 2903.6924  # J89: 3009.6155 -> 2903.6924
$ws.Cells.Item(89, 11).Value = 12169.1665  # K89: 13215.714 -> 12169.1665
$ws.Cells.Item(89, 12).Value = 14518.462  # L89: 15048.0775 -> 14518.462
$ws.Cells.Item(89, 13).Value = -6553.166499999999  # M89: -7599.714 -> -6553.166499999999
$ws.Cells.Item(89, 14).Value = -25750.462  # N89: -26280.0775 -> -25750.462
$ws.Cells.Item(94, 8).Value = 989  # H94: 983 -> 989
$ws.Cells.Item(94, 9).Value = 989  # I94: 983 -> 989
$ws.Cells.Item(94, 11).Value = 989  # K94: 983 -> 989
$ws.Cells.Item(94, 13).Value = -538  # M94: -532 -> -538
$ws.Cells.Item(98, 8).Value = 3621.7222  # H98: 3711.647 -> 3621.7222
$ws.Cells.Item(98, 9).Value = 3706.0588  # I98: 3711.647 -> 3706.0588
$ws.Cells.Item(98, 10).Value = 2188  # J98: 0 -> 2188
$ws.Cells.Item(98, 11).Value = 3706.0588  # K98: 3711.647 -> 3706.0588
$ws.Cells.Item(98, 12).Value = 2188  # L98: 0 -> 2188
$ws.Cells.Item(98, 13).Value = -2208.0588  # M98: -2213.647 -> -2208.0588
$ws.Cells.Item(98, 14).Value = -5184  # N98: (empty) -> -5184
$ws.Cells.Item(100, 8).Value = 9187.143  # H100: 11372 -> 9187.143
$ws.Cells.Item(100, 9).Value = 13707.5  # I100: 17760 -> 13707.5
$ws.Cells.Item(100, 10).Value = 3160  # J100: 1790 -> 3160
$ws.Cells.Item(100, 11).Value = 13707.5  # K100: 17760 -> 13707.5
$ws.Cells.Item(100, 12).Value = 3160  # L100: 1790 -> 3160
$ws.Cells.Item(100, 13).Value = -13166.5  # M100: -17219 -> -13166.5
$ws.Cells.Item(100, 14).Value = -4242  # N100: -2872 -> -4242
$ws.Cells.Item(103, 8).Value = 756.24445  # H103: 782 -> 756.24445
$ws.Cells.Item(103, 9).Value = 519.19354  # I103: 526.6667 -> 519.19354
$ws.Cells.Item(103, 10).Value = 1281.1428  # J103: 1371.2307 -> 1281.1428
$ws.Cells.Item(103, 11).Value = 1557.58062  # K103: 1580.0001 -> 1557.58062
$ws.Cells.Item(103, 12).Value = 3843.4284  # L103: 4113.6921 -> 3843.4284
$ws.Cells.Item(103, 13).Value = -971.58062  # M103: -994.0001 -> -971.58062
$ws.Cells.Item(103, 14).Value = -5015.428400000001  # N103: -5285.6921 -> -5015.428400000001
$ws.Cells.Item(113, 8).Value = 4792.857  # H113: 4801 -> 4792.857
$ws.Cells.Item(113, 10).Value = 5350.2  # J113: 5310.636 -> 5350.2
$ws.Cells.Item(113, 12).Value = 5350.2  # L113: 5310.636 -> 5350.2
$ws.Cells.Item(113, 14).Value = -11858.2  # N113: -11818.636 -> -11858.2
$ws.Cells.Item(116, 8).Value = 3000  # H116: 0 -> 3000
$ws.Cells.Item(116, 9).Value = 3000  # I116: 0 -> 3000
$ws.Cells.Item(116, 11).Value = 3000  # K116: 0 -> 3000
$ws.Cells.Item(116, 13).Value = 442  # M116: (empty) -> 442
$ws.Cells.Item(122, 8).Value = 3621.7222  # H122: 3711.647 -> 3621.7222
$ws.Cells.Item(122, 9).Value = 3706.0588  # I122: 3711.647 -> 3706.0588
$ws.Cells.Item(122, 10).Value = 2188  # J122: 0 -> 2188
$ws.Cells.Item(122, 11).Value = 11118.1764  # K122: 11134.941 -> 11118.1764
$ws.Cells.Item(122, 12).Value = 6564  # L122: 0 -> 6564
$ws.Cells.Item(122, 13).Value = -8668.1764  # M122: -8684.940999999999 -> -8668.1764
$ws.Cells.Item(122, 14).Value = -11464  # N122: (empty) -> -11464
$ws.Cells.Item(125, 8).Value = 1125.5  # H125: 1103.1 -> 1125.5
$ws.Cells.Item(125, 9).Value = 727  # I125: 644 -> 727
$ws.Cells.Item(125, 10).Value = 1364.6  # J125: 1562.2 -> 1364.6
$ws.Cells.Item(125, 11).Value = 6543  # K125: 5796 -> 6543
$ws.Cells.Item(125, 12).Value = 12281.4  # L125: 14059.8 -> 12281.4
$ws.Cells.Item(125, 13).Value = -4083  # M125: -3336 -> -4083
$ws.Cells.Item(125, 14).Value = -17201.4  # N125: -18979.8 -> -17201.4
$ws.Cells.Item(138, 8).Value = 2479.44  # H138: 2519.926 -> 2479.44
$ws.Cells.Item(138, 9).Value = 556.9231  # I138: 567.36 -> 556.9231
$ws.Cells.Item(138, 10).Value = 4562.1665  # J138: 4203.1724 -> 4562.1665
$ws.Cells.Item(138, 11).Value = 1670.7693  # K138: 1702.08 -> 1670.7693
$ws.Cells.Item(138, 12).Value = 13686.4995  # L138: 12609.5172 -> 13686.4995
$ws.Cells.Item(138, 13).Value = 3469.2307  # M138: 3437.92 -> 3469.2307
$ws.Cells.Item(138, 14).Value = -23966.4995  # N138: -22889.5172 -> -23966.4995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1417.6818  # H2: 1351.4348 -> 1417.6818
$ws.Cells.Item(2, 9).Value = 1746.9166  # I2: 1511.2142 -> 1746.9166
$ws.Cells.Item(2, 10).Value = 1022.6  # J2: 1102.8889 -> 1022.6
$ws.Cells.Item(2, 11).Value = 1746.9166  # K2: 1511.2142 -> 1746.9166
$ws.Cells.Item(2, 12).Value = 1022.6  # L2: 1102.8889 -> 1022.6
$ws.Cells.Item(2, 13).Value = -1633.9166  # M2: -1398.2142 -> -1633.9166
$ws.Cells.Item(2, 14).Value = -1248.6  # N2: -1328.8889 -> -1248.6
$ws.Cells.Item(7, 8).Value = 31384.615  # H7: 30000 -> 31384.615
$ws.Cells.Item(7, 10).Value = 31384.615  # J7: 30000 -> 31384.615
$ws.Cells.Item(7, 12).Value = 31384.615  # L7: 30000 -> 31384.615
$ws.Cells.Item(7, 14).Value = -31612.615  # N7: -30228 -> -31612.615
$ws.Cells.Item(32, 8).Value = 2044625.6  # H32: 1964473.4 -> 2044625.6
$ws.Cells.Item(32, 9).Value = 2044625.6  # I32: 1964473.4 -> 2044625.6
$ws.Cells.Item(32, 11).Value = 2044625.6  # K32: 1964473.4 -> 2044625.6
$ws.Cells.Item(32, 13).Value = -2044338.6  # M32: -1964186.4 -> -2044338.6
$ws.Cells.Item(43, 8).Value = 30381.834  # H43: 27569.4 -> 30381.834
$ws.Cells.Item(43, 10).Value = 31559  # J43: 28337.75 -> 31559
$ws.Cells.Item(43, 12).Value = 31559  # L43: 28337.75 -> 31559
$ws.Cells.Item(43, 14).Value = -32185  # N43: -28963.75 -> -32185
$ws.Cells.Item(45, 8).Value = 3362.182  # H45: 2612.9473 -> 3362.182
$ws.Cells.Item(45, 9).Value = 2366.1667  # I45: 1988.25 -> 2366.1667
$ws.Cells.Item(45, 10).Value = 4557.4  # J45: 3683.8572 -> 4557.4
$ws.Cells.Item(45, 11).Value = 2366.1667  # K45: 1988.25 -> 2366.1667
$ws.Cells.Item(45, 12).Value = 4557.4  # L45: 3683.8572 -> 4557.4
$ws.Cells.Item(45, 13).Value = -1989.1667  # M45: -1611.25 -> -1989.1667
$ws.Cells.Item(45, 14).Value = -5311.4  # N45: -4437.8572 -> -5311.4
$ws.Cells.Item(68, 8).Value = 30000  # H68: 100000 -> 30000
$ws.Cells.Item(68, 9).Value = 30000  # I68: 0 -> 30000
$ws.Cells.Item(68, 10).Value = 0  # J68: 100000 -> 0
$ws.Cells.Item(68, 11).Value = 30000  # K68: 0 -> 30000
$ws.Cells.Item(68, 12).Value = 0  # L68: 100000 -> 0
$ws.Cells.Item(68, 13).Value = -29189  # M68: (empty) -> -29189
$ws.Cells.Item(68, 14).ClearContents()  # N68: -101622 -> (cleared)
$ws.Cells.Item(71, 8).Value = 30000  # H71: 100000 -> 30000
$ws.Cells.Item(71, 9).Value = 30000  # I71: 0 -> 30000
$ws.Cells.Item(71, 10).Value = 0  # J71: 100000 -> 0
$ws.Cells.Item(71, 11).Value = 90000  # K71: 0 -> 90000
$ws.Cells.Item(71, 12).Value = 0  # L71: 300000 -> 0
$ws.Cells.Item(71, 13).Value = -85944  # M71: (empty) -> -85944
$ws.Cells.Item(71, 14).ClearContents()  # N71: -308112 -> (cleared)
$ws.Cells.Item(74, 8).Value = 1553.6451  # H74: 1574.2623 -> 1553.6451
$ws.Cells.Item(74, 9).Value = 943.3333  # I74: 960.3684 -> 943.3333
$ws.Cells.Item(74, 11).Value = 943.3333  # K74: 960.3684 -> 943.3333
$ws.Cells.Item(74, 13).Value = -69.33330000000001  # M74: -86.36839999999995 -> -69.33330000000001
$ws.Cells.Item(77, 8).Value = 1553.6451  # H77: 1574.2623 -> 1553.6451
$ws.Cells.Item(77, 9).Value = 943.3333  # I77: 960.3684 -> 943.3333
$ws.Cells.Item(77, 11).Value = 4716.6665  # K77: 4801.842 -> 4716.6665
$ws.Cells.Item(77, 13).Value = -348.6665000000003  # M77: -433.8419999999996 -> -348.6665000000003
$ws.Cells.Item(81, 8).Value = 99849.5  # H81: 99900 -> 99849.5
$ws.Cells.Item(81, 10).Value = 99849.5  # J81: 99900 -> 99849.5
$ws.Cells.Item(81, 12).Value = 99849.5  # L81: 99900 -> 99849.5
$ws.Cells.Item(81, 14).Value = -101845.5  # N81: -101896 -> -101845.5
$ws.Cells.Item(84, 8).Value = 99849.5  # H84: 99900 -> 99849.5
$ws.Cells.Item(84, 10).Value = 99849.5  # J84: 99900 -> 99849.5
$ws.Cells.Item(84, 12).Value = 299548.5  # L84: 299700 -> 299548.5
$ws.Cells.Item(84, 14).Value = -309532.5  # N84: -309684 -> -309532.5
$ws.Cells.Item(102, 8).Value = 1058.4546  # H102: 1081.9 -> 1058.4546
$ws.Cells.Item(102, 9).Value = 1053.2  # I102: 1078.6666 -> 1053.2
$ws.Cells.Item(102, 11).Value = 1053.2  # K102: 1078.6666 -> 1053.2
$ws.Cells.Item(102, 13).Value = 568.8  # M102: 543.3334 -> 568.8
$ws.Cells.Item(116, 8).Value = 1417.6818  # H116: 1351.4348 -> 1417.6818
$ws.Cells.Item(116, 9).Value = 1746.9166  # I116: 1511.2142 -> 1746.9166
$ws.Cells.Item(116, 10).Value = 1022.6  # J116: 1102.8889 -> 1022.6
$ws.Cells.Item(116, 11).Value = 1746.9166  # K116: 1511.2142 -> 1746.9166
$ws.Cells.Item(116, 12).Value = 1022.6  # L116: 1102.8889 -> 1022.6
$ws.Cells.Item(116, 13).Value = 547.0834  # M116: 782.7858000000001 -> 547.0834
$ws.Cells.Item(116, 14).Value = -5610.6  # N116: -5690.8889 -> -5610.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1417.6818  # H3: 1351.4348 -> 1417.6818
$ws.Cells.Item(3, 9).Value = 1746.9166  # I3: 1511.2142 -> 1746.9166
$ws.Cells.Item(3, 10).Value = 1022.6  # J3: 1102.8889 -> 1022.6
$ws.Cells.Item(3, 11).Value = 1746.9166  # K3: 1511.2142 -> 1746.9166
$ws.Cells.Item(3, 12).Value = 1022.6  # L3: 1102.8889 -> 1022.6
$ws.Cells.Item(3, 13).Value = -1632.9166  # M3: -1397.2142 -> -1632.9166
$ws.Cells.Item(3, 14).Value = -1250.6  # N3: -1330.8889 -> -1250.6
$ws.Cells.Item(50, 8).Value = 48885  # H50: 0 -> 48885
$ws.Cells.Item(50, 10).Value = 48885  # J50: 0 -> 48885
$ws.Cells.Item(50, 12).Value = 48885  # L50: 0 -> 48885
$ws.Cells.Item(50, 14).Value = -50033  # N50: (empty) -> -50033
$ws.Cells.Item(82, 8).Value = 6467.8184  # H82: 6839.6 -> 6467.8184
$ws.Cells.Item(82, 9).Value = 4486.3  # I82: 4679.222 -> 4486.3
$ws.Cells.Item(82, 11).Value = 4486.3  # K82: 4679.222 -> 4486.3
$ws.Cells.Item(82, 13).Value = -4103.3  # M82: -4296.222 -> -4103.3
$ws.Cells.Item(85, 8).Value = 6467.8184  # H85: 6839.6 -> 6467.8184
$ws.Cells.Item(85, 9).Value = 4486.3  # I85: 4679.222 -> 4486.3
$ws.Cells.Item(85, 11).Value = 4486.3  # K85: 4679.222 -> 4486.3
$ws.Cells.Item(85, 13).Value = -3160.3  # M85: -3353.222 -> -3160.3
$ws.Cells.Item(94, 8).Value = 7069.1  # H94: 7087.65 -> 7069.1
$ws.Cells.Item(94, 9).Value = 2038.1538  # I94: 2066.6924 -> 2038.1538
$ws.Cells.Item(94, 11).Value = 2038.1538  # K94: 2066.6924 -> 2038.1538
$ws.Cells.Item(94, 13).Value = -1587.1538  # M94: -1615.6924 -> -1587.1538
$ws.Cells.Item(99, 8).Value = 820.7  # H99: 857.55554 -> 820.7
$ws.Cells.Item(99, 9).Value = 901  # I99: 959.8570999999999 -> 901
$ws.Cells.Item(99, 11).Value = 901  # K99: 959.8570999999999 -> 901
$ws.Cells.Item(99, 13).Value = 597  # M99: 538.1429000000001 -> 597
$ws.Cells.Item(107, 8).Value = 2469.4583  # H107: 2314.1155 -> 2469.4583
$ws.Cells.Item(107, 9).Value = 1752.6875  # I107: 1607.9445 -> 1752.6875
$ws.Cells.Item(107, 11).Value = 1752.6875  # K107: 1607.9445 -> 1752.6875
$ws.Cells.Item(107, 13).Value = 167.3125  # M107: 312.0554999999999 -> 167.3125
$ws.Cells.Item(134, 8).Value = 48614200  # H134: 53033676 -> 48614200
$ws.Cells.Item(134, 9).Value = 27781062  # I134: 31253700 -> 27781062
$ws.Cells.Item(134, 11).Value = 83343186  # K134: 93761100 -> 83343186
$ws.Cells.Item(134, 13).Value = -83340651  # M134: -93758565 -> -83340651

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1499  # H16: 1238.7142 -> 1499
$ws.Cells.Item(16, 9).Value = 1499  # I16: 1200.8 -> 1499
$ws.Cells.Item(16, 10).Value = 0  # J16: 1333.5 -> 0
$ws.Cells.Item(16, 11).Value = 1499  # K16: 1200.8 -> 1499
$ws.Cells.Item(16, 12).Value = 0  # L16: 1333.5 -> 0
$ws.Cells.Item(16, 13).Value = -1212  # M16: -913.8 -> -1212
$ws.Cells.Item(16, 14).ClearContents()  # N16: -1907.5 -> (cleared)
$ws.Cells.Item(31, 8).Value = 1663.7142  # H31: 1718.0769 -> 1663.7142
$ws.Cells.Item(31, 9).Value = 1336.5  # I31: 1338.6 -> 1336.5
$ws.Cells.Item(31, 10).Value = 2481.75  # J31: 2983 -> 2481.75
$ws.Cells.Item(31, 11).Value = 1336.5  # K31: 1338.6 -> 1336.5
$ws.Cells.Item(31, 12).Value = 2481.75  # L31: 2983 -> 2481.75
$ws.Cells.Item(31, 13).Value = -1041.5  # M31: -1043.6 -> -1041.5
$ws.Cells.Item(31, 14).Value = -3071.75  # N31: -3573 -> -3071.75
$ws.Cells.Item(34, 8).Value = 1663.7142  # H34: 1718.0769 -> 1663.7142
$ws.Cells.Item(34, 9).Value = 1336.5  # I34: 1338.6 -> 1336.5
$ws.Cells.Item(34, 10).Value = 2481.75  # J34: 2983 -> 2481.75
$ws.Cells.Item(34, 11).Value = 1336.5  # K34: 1338.6 -> 1336.5
$ws.Cells.Item(34, 12).Value = 2481.75  # L34: 2983 -> 2481.75
$ws.Cells.Item(34, 13).Value = -1134.5  # M34: -1136.6 -> -1134.5
$ws.Cells.Item(34, 14).Value = -2885.75  # N34: -3387 -> -2885.75
$ws.Cells.Item(62, 8).Value = 46435.57  # H62: 46499.785 -> 46435.57
$ws.Cells.Item(62, 9).Value = 3853.6924  # I62: 3908.1667 -> 3853.6924
$ws.Cells.Item(62, 10).Value = 600000  # J62: 302049.5 -> 600000
$ws.Cells.Item(62, 11).Value = 3853.6924  # K62: 3908.1667 -> 3853.6924
$ws.Cells.Item(62, 12).Value = 600000  # L62: 302049.5 -> 600000
$ws.Cells.Item(62, 13).Value = -3229.6924  # M62: -3284.1667 -> -3229.6924
$ws.Cells.Item(62, 14).Value = -601248  # N62: -303297.5 -> -601248
$ws.Cells.Item(65, 8).Value = 46435.57  # H65: 46499.785 -> 46435.57
$ws.Cells.Item(65, 9).Value = 3853.6924  # I65: 3908.1667 -> 3853.6924
$ws.Cells.Item(65, 10).Value = 600000  # J65: 302049.5 -> 600000
$ws.Cells.Item(65, 11).Value = 19268.462  # K65: 19540.8335 -> 19268.462
$ws.Cells.Item(65, 12).Value = 3000000  # L65: 1510247.5 -> 3000000
$ws.Cells.Item(65, 13).Value = -16148.462  # M65: -16420.8335 -> -16148.462
$ws.Cells.Item(65, 14).Value = -3006240  # N65: -1516487.5 -> -3006240
$ws.Cells.Item(105, 8).Value = 1236.7142  # H105: 1236.9286 -> 1236.7142
$ws.Cells.Item(105, 9).Value = 763.7  # I105: 735 -> 763.7
$ws.Cells.Item(105, 10).Value = 2419.25  # J105: 3077.3333 -> 2419.25
$ws.Cells.Item(105, 11).Value = 763.7  # K105: 735 -> 763.7
$ws.Cells.Item(105, 12).Value = 2419.25  # L105: 3077.3333 -> 2419.25
$ws.Cells.Item(105, 13).Value = 983.3  # M105: 1012 -> 983.3
$ws.Cells.Item(105, 14).Value = -5913.25  # N105: -6571.3333 -> -5913.25
$ws.Cells.Item(113, 8).Value = 1499  # H113: 1238.7142 -> 1499
$ws.Cells.Item(113, 9).Value = 1499  # I113: 1200.8 -> 1499
$ws.Cells.Item(113, 10).Value = 0  # J113: 1333.5 -> 0
$ws.Cells.Item(113, 11).Value = 1499  # K113: 1200.8 -> 1499
$ws.Cells.Item(113, 12).Value = 0  # L113: 1333.5 -> 0
$ws.Cells.Item(113, 13).Value = 671  # M113: 969.2 -> 671
$ws.Cells.Item(113, 14).ClearContents()  # N113: -5673.5 -> (cleared)
$ws.Cells.Item(135, 8).Value = 80780  # H135: 0 -> 80780
$ws.Cells.Item(135, 10).Value = 80780  # J135: 0 -> 80780
$ws.Cells.Item(135, 12).Value = 80780  # L135: 0 -> 80780
$ws.Cells.Item(135, 14).Value = -90920  # N135: (empty) -> -90920

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 157681.42  # H4: 143409.92 -> 157681.42
$ws.Cells.Item(4, 9).Value = 20304  # I4: 323.9 -> 20304
$ws.Cells.Item(4, 11).Value = 60912  # K4: 971.6999999999999 -> 60912
$ws.Cells.Item(4, 13).Value = -60800  # M4: -859.6999999999999 -> -60800
$ws.Cells.Item(44, 8).Value = 3333  # H44: 922.5 -> 3333
$ws.Cells.Item(44, 9).Value = 0  # I44: 922.5 -> 0
$ws.Cells.Item(44, 10).Value = 3333  # J44: 0 -> 3333
$ws.Cells.Item(44, 11).Value = 0  # K44: 2767.5 -> 0
$ws.Cells.Item(44, 12).Value = 9999  # L44: 0 -> 9999
$ws.Cells.Item(44, 13).ClearContents()  # M44: -2369.5 -> (cleared)
$ws.Cells.Item(44, 14).Value = -10795  # N44: (empty) -> -10795
$ws.Cells.Item(69, 8).Value = 1000  # H69: 999.75 -> 1000
$ws.Cells.Item(69, 9).Value = 1000  # I69: 999.75 -> 1000
$ws.Cells.Item(69, 11).Value = 3000  # K69: 2999.25 -> 3000
$ws.Cells.Item(69, 13).Value = -2189  # M69: -2188.25 -> -2189
$ws.Cells.Item(72, 8).Value = 1000  # H72: 999.75 -> 1000
$ws.Cells.Item(72, 9).Value = 1000  # I72: 999.75 -> 1000
$ws.Cells.Item(72, 11).Value = 9000  # K72: 8997.75 -> 9000
$ws.Cells.Item(72, 13).Value = -4944  # M72: -4941.75 -> -4944
$ws.Cells.Item(74, 8).Value = 9999.166999999999  # H74: 10000 -> 9999.166999999999
$ws.Cells.Item(74, 10).Value = 9999.166999999999  # J74: 10000 -> 9999.166999999999
$ws.Cells.Item(74, 12).Value = 29997.501  # L74: 30000 -> 29997.501
$ws.Cells.Item(74, 14).Value = -32119.501  # N74: -32122 -> -32119.501
$ws.Cells.Item(77, 8).Value = 9999.166999999999  # H77: 10000 -> 9999.166999999999
$ws.Cells.Item(77, 10).Value = 9999.166999999999  # J77: 10000 -> 9999.166999999999
$ws.Cells.Item(77, 12).Value = 89992.503  # L77: 90000 -> 89992.503
$ws.Cells.Item(77, 14).Value = -100600.503  # N77: -100608 -> -100600.503
$ws.Cells.Item(81, 8).Value = 10998.333  # H81: 10999 -> 10998.333
$ws.Cells.Item(81, 10).Value = 10998.333  # J81: 10999 -> 10998.333
$ws.Cells.Item(81, 12).Value = 32994.999  # L81: 32997 -> 32994.999
$ws.Cells.Item(81, 14).Value = -35240.999  # N81: -35243 -> -35240.999
$ws.Cells.Item(84, 8).Value = 10998.333  # H84: 10999 -> 10998.333
$ws.Cells.Item(84, 10).Value = 10998.333  # J84: 10999 -> 10998.333
$ws.Cells.Item(84, 12).Value = 98984.997  # L84: 98991 -> 98984.997
$ws.Cells.Item(84, 14).Value = -110216.997  # N84: -110223 -> -110216.997
$ws.Cells.Item(107, 8).Value = 782.381  # H107: 730.34784 -> 782.381
$ws.Cells.Item(107, 9).Value = 353.57144  # I107: 334.125 -> 353.57144
$ws.Cells.Item(107, 10).Value = 996.7857  # J107: 941.6667 -> 996.7857
$ws.Cells.Item(107, 11).Value = 1060.71432  # K107: 1002.375 -> 1060.71432
$ws.Cells.Item(107, 12).Value = 2990.3571  # L107: 2825.0001 -> 2990.3571
$ws.Cells.Item(107, 13).Value = 859.28568  # M107: 917.625 -> 859.28568
$ws.Cells.Item(107, 14).Value = -6830.3571  # N107: -6665.0001 -> -6830.3571
$ws.Cells.Item(131, 8).Value = 254082.52  # H131: 265074.66 -> 254082.52
$ws.Cells.Item(131, 9).Value = 996.3333  # I131: 997.6 -> 996.3333
$ws.Cells.Item(131, 10).Value = 277090.34  # J131: 285705.7 -> 277090.34
$ws.Cells.Item(131, 11).Value = 2988.9999  # K131: 2992.8 -> 2988.9999
$ws.Cells.Item(131, 12).Value = 831271.02  # L131: 857117.1000000001 -> 831271.02
$ws.Cells.Item(131, 13).Value = 2051.0001  # M131: 2047.2 -> 2051.0001
$ws.Cells.Item(131, 14).Value = -841351.02  # N131: -867197.1000000001 -> -841351.02
$ws.Cells.Item(132, 8).Value = 12958  # H132: 14649.6 -> 12958
$ws.Cells.Item(132, 10).Value = 32125  # J132: 59750 -> 32125
$ws.Cells.Item(132, 12).Value = 289125  # L132: 537750 -> 289125
$ws.Cells.Item(132, 14).Value = -294185  # N132: -542810 -> -294185

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 11429.6  # H22: 13599.5 -> 11429.6
$ws.Cells.Item(22, 9).Value = 3749.3333  # I22: 4249 -> 3749.3333
$ws.Cells.Item(22, 11).Value = 3749.3333  # K22: 4249 -> 3749.3333
$ws.Cells.Item(22, 13).Value = -3220.3333  # M22: -3720 -> -3220.3333
$ws.Cells.Item(74, 8).Value = 100000  # H74: 0 -> 100000
$ws.Cells.Item(74, 10).Value = 100000  # J74: 0 -> 100000
$ws.Cells.Item(74, 12).Value = 100000  # L74: 0 -> 100000
$ws.Cells.Item(74, 14).Value = -101872  # N74: (empty) -> -101872
$ws.Cells.Item(77, 8).Value = 100000  # H77: 0 -> 100000
$ws.Cells.Item(77, 10).Value = 100000  # J77: 0 -> 100000
$ws.Cells.Item(77, 12).Value = 300000  # L77: 0 -> 300000
$ws.Cells.Item(77, 14).Value = -309360  # N77: (empty) -> -309360
$ws.Cells.Item(123, 8).Value = 65000  # H123: 0 -> 65000
$ws.Cells.Item(123, 10).Value = 65000  # J123: 0 -> 65000
$ws.Cells.Item(123, 12).Value = 65000  # L123: 0 -> 65000
$ws.Cells.Item(123, 14).Value = -69900  # N123: (empty) -> -69900
$ws.Cells.Item(124, 8).Value = 57597.125  # H124: 59999.855 -> 57597.125
$ws.Cells.Item(124, 10).Value = 57597.125  # J124: 59999.855 -> 57597.125
$ws.Cells.Item(124, 12).Value = 57597.125  # L124: 59999.855 -> 57597.125
$ws.Cells.Item(124, 14).Value = -67417.125  # N124: -69819.85500000001 -> -67417.125
$ws.Cells.Item(132, 8).Value = 1947.9445  # H132: 1997.1765 -> 1947.9445
$ws.Cells.Item(132, 9).Value = 1698.1428  # I132: 1743.3077 -> 1698.1428
$ws.Cells.Item(132, 11).Value = 5094.428400000001  # K132: 5229.9231 -> 5094.428400000001
$ws.Cells.Item(132, 13).Value = -2564.428400000001  # M132: -2699.9231 -> -2564.428400000001
$ws.Cells.Item(136, 8).Value = 33169.08  # H136: 33709.08 -> 33169.08
$ws.Cells.Item(136, 10).Value = 33169.08  # J136: 33709.08 -> 33169.08
$ws.Cells.Item(136, 12).Value = 99507.24000000001  # L136: 101127.24 -> 99507.24000000001
$ws.Cells.Item(136, 14).Value = -104607.24  # N136: -106227.24 -> -104607.24

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3961.8333  # H7: 4356.2 -> 3961.8333
$ws.Cells.Item(7, 9).Value = 4353.2  # I7: 4356.2 -> 4353.2
$ws.Cells.Item(7, 10).Value = 2005  # J7: 0 -> 2005
$ws.Cells.Item(7, 11).Value = 4353.2  # K7: 4356.2 -> 4353.2
$ws.Cells.Item(7, 12).Value = 2005  # L7: 0 -> 2005
$ws.Cells.Item(7, 13).Value = -4241.2  # M7: -4244.2 -> -4241.2
$ws.Cells.Item(7, 14).Value = -2229  # N7: (empty) -> -2229
$ws.Cells.Item(14, 8).Value = 23329.834  # H14: 45200 -> 23329.834
$ws.Cells.Item(14, 10).Value = 24195.8  # J14: 51750 -> 24195.8
$ws.Cells.Item(14, 12).Value = 24195.8  # L14: 51750 -> 24195.8
$ws.Cells.Item(14, 14).Value = -24539.8  # N14: -52094 -> -24539.8
$ws.Cells.Item(22, 8).Value = 2176.8  # H22: 1697.8182 -> 2176.8
$ws.Cells.Item(22, 9).Value = 2176.8  # I22: 1910 -> 2176.8
$ws.Cells.Item(22, 10).Value = 0  # J22: 1132 -> 0
$ws.Cells.Item(22, 11).Value = 2176.8  # K22: 1910 -> 2176.8
$ws.Cells.Item(22, 12).Value = 0  # L22: 1132 -> 0
$ws.Cells.Item(22, 13).Value = -1881.8  # M22: -1615 -> -1881.8
$ws.Cells.Item(22, 14).ClearContents()  # N22: -1722 -> (cleared)
$ws.Cells.Item(27, 8).Value = 2176.8  # H27: 1697.8182 -> 2176.8
$ws.Cells.Item(27, 9).Value = 2176.8  # I27: 1910 -> 2176.8
$ws.Cells.Item(27, 10).Value = 0  # J27: 1132 -> 0
$ws.Cells.Item(27, 11).Value = 2176.8  # K27: 1910 -> 2176.8
$ws.Cells.Item(27, 12).Value = 0  # L27: 1132 -> 0
$ws.Cells.Item(27, 13).Value = -2069.8  # M27: -1803 -> -2069.8
$ws.Cells.Item(27, 14).ClearContents()  # N27: -1346 -> (cleared)
$ws.Cells.Item(46, 8).Value = 3156.875  # H46: 3913.0527 -> 3156.875
$ws.Cells.Item(46, 9).Value = 840  # I46: 1333.6666 -> 840
$ws.Cells.Item(46, 10).Value = 3766.5789  # J46: 4396.6875 -> 3766.5789
$ws.Cells.Item(46, 11).Value = 840  # K46: 1333.6666 -> 840
$ws.Cells.Item(46, 12).Value = 3766.5789  # L46: 4396.6875 -> 3766.5789
$ws.Cells.Item(46, 13).Value = -652  # M46: -1145.6666 -> -652
$ws.Cells.Item(46, 14).Value = -4142.5789  # N46: -4772.6875 -> -4142.5789
$ws.Cells.Item(61, 8).Value = 2118.25  # H61: 2078.524 -> 2118.25
$ws.Cells.Item(61, 9).Value = 1485.5294  # I61: 1474.3334 -> 1485.5294
$ws.Cells.Item(61, 11).Value = 1485.5294  # K61: 1474.3334 -> 1485.5294
$ws.Cells.Item(61, 13).Value = -1283.5294  # M61: -1272.3334 -> -1283.5294
$ws.Cells.Item(68, 8).Value = 2614.4  # H68: 2558.5454 -> 2614.4
$ws.Cells.Item(68, 9).Value = 2738.3333  # I68: 2664.5 -> 2738.3333
$ws.Cells.Item(68, 11).Value = 2738.3333  # K68: 2664.5 -> 2738.3333
$ws.Cells.Item(68, 13).Value = -1989.3333  # M68: -1915.5 -> -1989.3333
$ws.Cells.Item(71, 8).Value = 2614.4  # H71: 2558.5454 -> 2614.4
$ws.Cells.Item(71, 9).Value = 2738.3333  # I71: 2664.5 -> 2738.3333
$ws.Cells.Item(71, 11).Value = 13691.6665  # K71: 13322.5 -> 13691.6665
$ws.Cells.Item(71, 13).Value = -9947.666499999999  # M71: -9578.5 -> -9947.666499999999
$ws.Cells.Item(93, 8).Value = 15038.207  # H93: 16094.963 -> 15038.207
$ws.Cells.Item(93, 9).Value = 1102.8422  # I93: 1128.3334 -> 1102.8422
$ws.Cells.Item(93, 10).Value = 41515.4  # J93: 46028.223 -> 41515.4
$ws.Cells.Item(93, 11).Value = 1102.8422  # K93: 1128.3334 -> 1102.8422
$ws.Cells.Item(93, 12).Value = 41515.4  # L93: 46028.223 -> 41515.4
$ws.Cells.Item(93, 13).Value = 145.1578  # M93: 119.6666 -> 145.1578
$ws.Cells.Item(93, 14).Value = -44011.4  # N93: -48524.223 -> -44011.4
$ws.Cells.Item(108, 8).Value = 52997.5  # H108: 46816.5 -> 52997.5
$ws.Cells.Item(108, 10).Value = 52997.5  # J108: 46816.5 -> 52997.5
$ws.Cells.Item(108, 12).Value = 52997.5  # L108: 46816.5 -> 52997.5
$ws.Cells.Item(108, 14).Value = -60677.5  # N108: -54496.5 -> -60677.5
$ws.Cells.Item(113, 8).Value = 2118.25  # H113: 2078.524 -> 2118.25
$ws.Cells.Item(113, 9).Value = 1485.5294  # I113: 1474.3334 -> 1485.5294
$ws.Cells.Item(113, 11).Value = 1485.5294  # K113: 1474.3334 -> 1485.5294
$ws.Cells.Item(113, 13).Value = 684.4706000000001  # M113: 695.6666 -> 684.4706000000001
$ws.Cells.Item(122, 8).Value = 3136  # H122: 3344 -> 3136
$ws.Cells.Item(122, 9).Value = 2991.7  # I122: 3276.111 -> 2991.7
$ws.Cells.Item(122, 11).Value = 8975.099999999999  # K122: 9828.332999999999 -> 8975.099999999999
$ws.Cells.Item(122, 13).Value = -6525.099999999999  # M122: -7378.332999999999 -> -6525.099999999999
$ws.Cells.Item(126, 8).Value = 3961.8333  # H126: 4356.2 -> 3961.8333
$ws.Cells.Item(126, 9).Value = 4353.2  # I126: 4356.2 -> 4353.2
$ws.Cells.Item(126, 10).Value = 2005  # J126: 0 -> 2005
$ws.Cells.Item(126, 11).Value = 13059.6  # K126: 13068.6 -> 13059.6
$ws.Cells.Item(126, 12).Value = 6015  # L126: 0 -> 6015
$ws.Cells.Item(126, 13).Value = -10589.6  # M126: -10598.6 -> -10589.6
$ws.Cells.Item(126, 14).Value = -10955  # N126: (empty) -> -10955
$ws.Cells.Item(132, 8).Value = 6152.154  # H132: 5255.1875 -> 6152.154
$ws.Cells.Item(132, 9).Value = 3533  # I132: 2991.8333 -> 3533
$ws.Cells.Item(132, 11).Value = 10599  # K132: 8975.499899999999 -> 10599
$ws.Cells.Item(132, 13).Value = -8069  # M132: -6445.499899999999 -> -8069
$ws.Cells.Item(133, 8).Value = 97000  # H133: 0 -> 97000
$ws.Cells.Item(133, 10).Value = 97000  # J133: 0 -> 97000
$ws.Cells.Item(133, 12).Value = 97000  # L133: 0 -> 97000
$ws.Cells.Item(133, 14).Value = -102060  # N133: (empty) -> -102060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(43, 8).Value = 10000000  # H43: 0 -> 10000000
$ws.Cells.Item(43, 9).Value = 10000000  # I43: 0 -> 10000000
$ws.Cells.Item(43, 11).Value = 10000000  # K43: 0 -> 10000000
$ws.Cells.Item(43, 13).Value = -9999851  # M43: (empty) -> -9999851
$ws.Cells.Item(54, 8).Value = 20692.334  # H54: 21000 -> 20692.334
$ws.Cells.Item(54, 10).Value = 20692.334  # J54: 21000 -> 20692.334
$ws.Cells.Item(54, 12).Value = 20692.334  # L54: 21000 -> 20692.334
$ws.Cells.Item(54, 14).Value = -21732.334  # N54: -22040 -> -21732.334
$ws.Cells.Item(113, 8).Value = 1553.95  # H113: 1625.5264 -> 1553.95
$ws.Cells.Item(113, 9).Value = 611.7857  # I113: 643.9231 -> 611.7857
$ws.Cells.Item(113, 11).Value = 1835.3571  # K113: 1931.7693 -> 1835.3571
$ws.Cells.Item(113, 13).Value = 334.6428999999998  # M113: 238.2307000000001 -> 334.6428999999998
$ws.Cells.Item(126, 8).Value = 4456.7144  # H126: 4027 -> 4456.7144
$ws.Cells.Item(126, 9).Value = 6598.3335  # I126: 4390 -> 6598.3335
$ws.Cells.Item(126, 10).Value = 2850.5  # J126: 3301 -> 2850.5
$ws.Cells.Item(126, 11).Value = 19795.0005  # K126: 13170 -> 19795.0005
$ws.Cells.Item(126, 12).Value = 8551.5  # L126: 9903 -> 8551.5
$ws.Cells.Item(126, 13).Value = -17325.0005  # M126: -10700 -> -17325.0005
$ws.Cells.Item(126, 14).Value = -13491.5  # N126: -14843 -> -13491.5
$ws.Cells.Item(132, 8).Value = 3426.1304  # H132: 3490.2856 -> 3426.1304
$ws.Cells.Item(132, 9).Value = 3371.0476  # I132: 3464.55 -> 3371.0476
$ws.Cells.Item(132, 10).Value = 4004.5  # J132: 4005 -> 4004.5
$ws.Cells.Item(132, 11).Value = 10113.1428  # K132: 10393.65 -> 10113.1428
$ws.Cells.Item(132, 12).Value = 12013.5  # L132: 12015 -> 12013.5
$ws.Cells.Item(132, 13).Value = -7583.1428  # M132: -7863.650000000001 -> -7583.1428
$ws.Cells.Item(132, 14).Value = -17073.5  # N132: -17075 -> -17073.5
$ws.Cells.Item(136, 8).Value = 1961.1666  # H136: 2118.875 -> 1961.1666
$ws.Cells.Item(136, 9).Value = 1961.1666  # I136: 2118.875 -> 1961.1666
$ws.Cells.Item(136, 11).Value = 5883.4998  # K136: 6356.625 -> 5883.4998
$ws.Cells.Item(136, 13).Value = -3333.4998  # M136: -3806.625 -> -3333.4998

